$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.888950333333333
$ws.Cells.Item(2, 8).Value = 8.666850999999998
$ws.Cells.Item(2, 9).Value = 0.014845006111042317
$ws.Cells.Item(2, 10).Value = 0.014845006111042313
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.01339666666666667
$ws.Cells.Item(2, 14).Value = 0.04019000000000001
$ws.Cells.Item(2, 15).Value = 0.08393217762128818
$ws.Cells.Item(2, 16).Value = 0.0839321776212882
$ws.Cells.Item(2, 17).Value = 0.038702304632222226
$ws.Cells.Item(2, 18).Value = 0.34832074169
$ws.Cells.Item(2, 19).Value = 0.0012459736897011122
$ws.Cells.Item(2, 20).Value = 0.0012459736897011122

# Row 3
$ws.Cells.Item(3, 7).Value = 2.888950333333333
$ws.Cells.Item(3, 8).Value = 8.666850999999998
$ws.Cells.Item(3, 9).Value = 0.014845006111042317
$ws.Cells.Item(3, 10).Value = 0.014845006111042313
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.1260863333333333
$ws.Cells.Item(3, 14).Value = 0.3782589999999999
$ws.Cells.Item(3, 15).Value = 0.7899502755623496
$ws.Cells.Item(3, 16).Value = 0.7899502755623498
$ws.Cells.Item(3, 17).Value = 0.36425715471211095
$ws.Cells.Item(3, 18).Value = 3.278314392408998
$ws.Cells.Item(3, 19).Value = 0.011726816668142643
$ws.Cells.Item(3, 20).Value = 0.011726816668142641

# Row 4
$ws.Cells.Item(4, 7).Value = 2.888950333333333
$ws.Cells.Item(4, 8).Value = 8.666850999999998
$ws.Cells.Item(4, 9).Value = 0.014845006111042317
$ws.Cells.Item(4, 10).Value = 0.014845006111042313
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02013
$ws.Cells.Item(4, 14).Value = 0.06039
$ws.Cells.Item(4, 15).Value = 0.12611754681636209
$ws.Cells.Item(4, 16).Value = 0.12611754681636209
$ws.Cells.Item(4, 17).Value = 0.058154570209999985
$ws.Cells.Item(4, 18).Value = 0.5233911318899999
$ws.Cells.Item(4, 19).Value = 0.0018722157531985606
$ws.Cells.Item(4, 20).Value = 0.0018722157531985602

# Row 5
$ws.Cells.Item(5, 7).Value = 12.56197866666667
$ws.Cells.Item(5, 8).Value = 37.68593600000001
$ws.Cells.Item(5, 9).Value = 0.06455031362836974
$ws.Cells.Item(5, 10).Value = 0.06455031362836973
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01339666666666667
$ws.Cells.Item(5, 14).Value = 0.04019000000000001
$ws.Cells.Item(5, 15).Value = 0.08393217762128818
$ws.Cells.Item(5, 16).Value = 0.0839321776212882
$ws.Cells.Item(5, 17).Value = 0.1682886408711112
$ws.Cells.Item(5, 18).Value = 1.5145977678400009
$ws.Cells.Item(5, 19).Value = 0.0054178483889661885
$ws.Cells.Item(5, 20).Value = 0.0054178483889661885

# Row 6
$ws.Cells.Item(6, 7).Value = 12.56197866666667
$ws.Cells.Item(6, 8).Value = 37.68593600000001
$ws.Cells.Item(6, 9).Value = 0.06455031362836974
$ws.Cells.Item(6, 10).Value = 0.06455031362836973
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1260863333333333
$ws.Cells.Item(6, 14).Value = 0.3782589999999999
$ws.Cells.Item(6, 15).Value = 0.7899502755623496
$ws.Cells.Item(6, 16).Value = 0.7899502755623498
$ws.Cells.Item(6, 17).Value = 1.5838938294915557
$ws.Cells.Item(6, 18).Value = 14.255044465424001
$ws.Cells.Item(6, 19).Value = 0.05099153803836677
$ws.Cells.Item(6, 20).Value = 0.05099153803836677

# Row 7
$ws.Cells.Item(7, 7).Value = 12.56197866666667
$ws.Cells.Item(7, 8).Value = 37.68593600000001
$ws.Cells.Item(7, 9).Value = 0.06455031362836974
$ws.Cells.Item(7, 10).Value = 0.06455031362836973
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02013
$ws.Cells.Item(7, 14).Value = 0.06039
$ws.Cells.Item(7, 15).Value = 0.12611754681636209
$ws.Cells.Item(7, 16).Value = 0.12611754681636209
$ws.Cells.Item(7, 17).Value = 0.25287263056000003
$ws.Cells.Item(7, 18).Value = 2.275853675040001
$ws.Cells.Item(7, 19).Value = 0.008140927201036777
$ws.Cells.Item(7, 20).Value = 0.008140927201036775

# Row 8
$ws.Cells.Item(8, 7).Value = 94.40225766666667
$ws.Cells.Item(8, 8).Value = 283.206773
$ws.Cells.Item(8, 9).Value = 0.4850904066394559
$ws.Cells.Item(8, 10).Value = 0.48509040663945585
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.01339666666666667
$ws.Cells.Item(8, 14).Value = 0.04019000000000001
$ws.Cells.Item(8, 15).Value = 0.08393217762128818
$ws.Cells.Item(8, 16).Value = 0.0839321776212882
$ws.Cells.Item(8, 17).Value = 1.2646755785411115
$ws.Cells.Item(8, 18).Value = 11.382080206870002
$ws.Cells.Item(8, 19).Value = 0.04071469417244573
$ws.Cells.Item(8, 20).Value = 0.04071469417244573

# Row 9
$ws.Cells.Item(9, 7).Value = 94.40225766666667
$ws.Cells.Item(9, 8).Value = 283.206773
$ws.Cells.Item(9, 9).Value = 0.4850904066394559
$ws.Cells.Item(9, 10).Value = 0.48509040663945585
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.1260863333333333
$ws.Cells.Item(9, 14).Value = 0.3782589999999999
$ws.Cells.Item(9, 15).Value = 0.7899502755623496
$ws.Cells.Item(9, 16).Value = 0.7899502755623498
$ws.Cells.Item(9, 17).Value = 11.902834527578554
$ws.Cells.Item(9, 18).Value = 107.12551074820698
$ws.Cells.Item(9, 19).Value = 0.38319730039749045
$ws.Cells.Item(9, 20).Value = 0.38319730039749045

# Row 10
$ws.Cells.Item(10, 7).Value = 94.40225766666667
$ws.Cells.Item(10, 8).Value = 283.206773
$ws.Cells.Item(10, 9).Value = 0.4850904066394559
$ws.Cells.Item(10, 10).Value = 0.48509040663945585
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.02013
$ws.Cells.Item(10, 14).Value = 0.06039
$ws.Cells.Item(10, 15).Value = 0.12611754681636209
$ws.Cells.Item(10, 16).Value = 0.12611754681636209
$ws.Cells.Item(10, 17).Value = 1.90031744683
$ws.Cells.Item(10, 18).Value = 17.10285702147
$ws.Cells.Item(10, 19).Value = 0.061178412069519704
$ws.Cells.Item(10, 20).Value = 0.0611784120695197

# Row 11
$ws.Cells.Item(11, 7).Value = 0.421979
$ws.Cells.Item(11, 8).Value = 1.265937
$ws.Cells.Item(11, 9).Value = 0.0021683587846606086
$ws.Cells.Item(11, 10).Value = 0.0021683587846606086
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.01339666666666667
$ws.Cells.Item(11, 14).Value = 0.04019000000000001
$ws.Cells.Item(11, 15).Value = 0.08393217762128818
$ws.Cells.Item(11, 16).Value = 0.0839321776212882
$ws.Cells.Item(11, 17).Value = 0.0056531120033333345
$ws.Cells.Item(11, 18).Value = 0.05087800803000002
$ws.Cells.Item(11, 19).Value = 0.00018199507466081478
$ws.Cells.Item(11, 20).Value = 0.0001819950746608148

# Row 12
$ws.Cells.Item(12, 7).Value = 0.421979
$ws.Cells.Item(12, 8).Value = 1.265937
$ws.Cells.Item(12, 9).Value = 0.0021683587846606086
$ws.Cells.Item(12, 10).Value = 0.0021683587846606086
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.1260863333333333
$ws.Cells.Item(12, 14).Value = 0.3782589999999999
$ws.Cells.Item(12, 15).Value = 0.7899502755623496
$ws.Cells.Item(12, 16).Value = 0.7899502755623498
$ws.Cells.Item(12, 17).Value = 0.05320578485366665
$ws.Cells.Item(12, 18).Value = 0.47885206368299993
$ws.Cells.Item(12, 19).Value = 0.0017128956194606894
$ws.Cells.Item(12, 20).Value = 0.0017128956194606896

# Row 13
$ws.Cells.Item(13, 7).Value = 0.421979
$ws.Cells.Item(13, 8).Value = 1.265937
$ws.Cells.Item(13, 9).Value = 0.0021683587846606086
$ws.Cells.Item(13, 10).Value = 0.0021683587846606086
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.02013
$ws.Cells.Item(13, 14).Value = 0.06039
$ws.Cells.Item(13, 15).Value = 0.12611754681636209
$ws.Cells.Item(13, 16).Value = 0.12611754681636209
$ws.Cells.Item(13, 17).Value = 0.008494437269999999
$ws.Cells.Item(13, 18).Value = 0.07644993543
$ws.Cells.Item(13, 19).Value = 0.0002734680905391043
$ws.Cells.Item(13, 20).Value = 0.0002734680905391043

# Row 14
$ws.Cells.Item(14, 7).Value = 84.33238866666666
$ws.Cells.Item(14, 8).Value = 252.997166
$ws.Cells.Item(14, 9).Value = 0.4333459148364715
$ws.Cells.Item(14, 10).Value = 0.43334591483647145
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.01339666666666667
$ws.Cells.Item(14, 14).Value = 0.04019000000000001
$ws.Cells.Item(14, 15).Value = 0.08393217762128818
$ws.Cells.Item(14, 16).Value = 0.0839321776212882
$ws.Cells.Item(14, 17).Value = 1.1297729001711112
$ws.Cells.Item(14, 18).Value = 10.167956101540002
$ws.Cells.Item(14, 19).Value = 0.036371666295514346
$ws.Cells.Item(14, 20).Value = 0.03637166629551435

# Row 15
$ws.Cells.Item(15, 7).Value = 84.33238866666666
$ws.Cells.Item(15, 8).Value = 252.997166
$ws.Cells.Item(15, 9).Value = 0.4333459148364715
$ws.Cells.Item(15, 10).Value = 0.43334591483647145
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.1260863333333333
$ws.Cells.Item(15, 14).Value = 0.3782589999999999
$ws.Cells.Item(15, 15).Value = 0.7899502755623496
$ws.Cells.Item(15, 16).Value = 0.7899502755623498
$ws.Cells.Item(15, 17).Value = 10.633161668221552
$ws.Cells.Item(15, 18).Value = 95.69845501399398
$ws.Cells.Item(15, 19).Value = 0.3423217248388892
$ws.Cells.Item(15, 20).Value = 0.3423217248388892

# Row 16
$ws.Cells.Item(16, 7).Value = 84.33238866666666
$ws.Cells.Item(16, 8).Value = 252.997166
$ws.Cells.Item(16, 9).Value = 0.4333459148364715
$ws.Cells.Item(16, 10).Value = 0.43334591483647145
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.02013
$ws.Cells.Item(16, 14).Value = 0.06039
$ws.Cells.Item(16, 15).Value = 0.12611754681636209
$ws.Cells.Item(16, 16).Value = 0.12611754681636209
$ws.Cells.Item(16, 17).Value = 1.6976109838599998
$ws.Cells.Item(16, 18).Value = 15.278498854739999
$ws.Cells.Item(16, 19).Value = 0.054652523702067954
$ws.Cells.Item(16, 20).Value = 0.05465252370206795

